$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet 展览 (rows 3-14) ---
$ws1.Range("F3").Value = 52
$ws1.Range("F4").Value = 1455
$ws1.Range("F5").Value = 336
$ws1.Range("F6").Value = 1052
$ws1.Range("F7").Value = 10857
$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = '2024.04.06'
$ws1.Range("C8").Value = '苏州·第一届寒假动漫展宅舞比赛-CF01'
$ws1.Range("D8").Value = '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店'
$ws1.Range("E8").Value = '2024.04.06 10:00-04.06 16:00'
$ws1.Range("F8").Value = 85
$ws1.Range("G8").Value = 49
$ws1.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=80528'
$ws1.Range("I8").Value = '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'
$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = '2024.04.13'
$ws1.Range("C9").Value = '苏州·X-party 国漫游戏嘉年华03'
$ws1.Range("D9").Value = '秋枫街与开平路交叉口西南角 爱琴海购物中心'
$ws1.Range("E9").Value = '2024.04.13 10:00-04.14 17:00'
$ws1.Range("F9").Value = 29
$ws1.Range("G9").Value = 48
$ws1.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=82042'
$ws1.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202402/WaQk4nUt1708679999084.jpeg'
$ws1.Range("F10").Value = 303
$ws1.Range("G10").Value = '已停售'
$ws1.Range("F11").Value = 1054
$ws1.Range("F12").Value = 732
$ws1.Range("F13").Value = 12154
$ws1.Range("F14").Value = 12629

# --- Sheet 全部类型 (rows 4-15, offset +1 vs 展览) ---
$ws4.Range("F4").Value = 52
$ws4.Range("F5").Value = 1455
$ws4.Range("F6").Value = 336
$ws4.Range("F7").Value = 1052
$ws4.Range("F8").Value = 10857
$ws4.Range("B9").NumberFormat = "@"
$ws4.Range("B9").Value = '2024.04.06'
$ws4.Range("C9").Value = '苏州·第一届寒假动漫展宅舞比赛-CF01'
$ws4.Range("D9").Value = '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店'
$ws4.Range("E9").Value = '2024.04.06 10:00-04.06 16:00'
$ws4.Range("F9").Value = 85
$ws4.Range("G9").Value = 49
$ws4.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=80528'
$ws4.Range("I9").Value = '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'
$ws4.Range("B10").NumberFormat = "@"
$ws4.Range("B10").Value = '2024.04.13'
$ws4.Range("C10").Value = '苏州·X-party 国漫游戏嘉年华03'
$ws4.Range("D10").Value = '秋枫街与开平路交叉口西南角 爱琴海购物中心'
$ws4.Range("E10").Value = '2024.04.13 10:00-04.14 17:00'
$ws4.Range("F10").Value = 29
$ws4.Range("G10").Value = 48
$ws4.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=82042'
$ws4.Range("I10").Value = '//i1.hdslb.com/bfs/openplatform/202402/WaQk4nUt1708679999084.jpeg'
$ws4.Range("F11").Value = 303
$ws4.Range("G11").Value = '已停售'
$ws4.Range("F12").Value = 1054
$ws4.Range("F13").Value = 732
$ws4.Range("F14").Value = 12154
$ws4.Range("F15").Value = 12629
